$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I14").Value = "sv"
$ws.Range("J14").Value = "Statement-opinion"
$ws.Range("I20").Value = "aa"
$ws.Range("J20").Value = "Agree/Accept"
$ws.Range("I21").Value = "sd"
$ws.Range("J21").Value = "Statement-non-opinion"
$ws.Range("I25").Value = "sd"
$ws.Range("J25").Value = "Statement-non-opinion"
$ws.Range("I45").Value = "sd"
$ws.Range("J45").Value = "Statement-non-opinion"
$ws.Range("I53").Value = "aa"
$ws.Range("J53").Value = "Agree/Accept"
$ws.Range("I64").Value = "sd"
$ws.Range("J64").Value = "Statement-non-opinion"
$ws.Range("I89").Value = "sd"
$ws.Range("J89").Value = "Statement-non-opinion"
$ws.Range("I92").Value = "sd"
$ws.Range("J92").Value = "Statement-non-opinion"
$ws.Range("I93").Value = "sd"
$ws.Range("J93").Value = "Statement-non-opinion"
$ws.Range("I101").Value = "sd"
$ws.Range("J101").Value = "Statement-non-opinion"
$ws.Range("I102").Value = "sd"
$ws.Range("J102").Value = "Statement-non-opinion"
$ws.Range("I105").Value = "sd"
$ws.Range("J105").Value = "Statement-non-opinion"
$ws.Range("I106").Value = "sd"
$ws.Range("J106").Value = "Statement-non-opinion"
$ws.Range("I120").Value = "sd"
$ws.Range("J120").Value = "Statement-non-opinion"
$ws.Range("I122").Value = "sd"
$ws.Range("J122").Value = "Statement-non-opinion"
$ws.Range("I130").Value = "sd"
$ws.Range("J130").Value = "Statement-non-opinion"
$ws.Range("I131").Value = "sd"
$ws.Range("J131").Value = "Statement-non-opinion"
$ws.Range("I132").Value = "sd"
$ws.Range("J132").Value = "Statement-non-opinion"
$ws.Range("I133").Value = "sd"
$ws.Range("J133").Value = "Statement-non-opinion"
$ws.Range("I134").Value = "sd"
$ws.Range("J134").Value = "Statement-non-opinion"
$ws.Range("I135").Value = "sv"
$ws.Range("J135").Value = "Statement-opinion"
$ws.Range("I144").Value = "sv"
$ws.Range("J144").Value = "Statement-opinion"
$ws.Range("I147").Value = "sv"
$ws.Range("J147").Value = "Statement-opinion"
$ws.Range("I150").Value = "b"
$ws.Range("J150").Value = "Acknowledge (Backchannel)"
$ws.Range("I154").Value = "sv"
$ws.Range("J154").Value = "Statement-opinion"
